$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I13").Value = 'b'
$ws.Range("J13").Value = 'Acknowledge (Backchannel)'
$ws.Range("I50").Value = 'b'
$ws.Range("J50").Value = 'Acknowledge (Backchannel)'
$ws.Range("I56").Value = 'sv'
$ws.Range("J56").Value = 'Statement-opinion'
$ws.Range("I57").Value = 'b'
$ws.Range("J57").Value = 'Acknowledge (Backchannel)'
$ws.Range("I60").Value = 'b'
$ws.Range("J60").Value = 'Acknowledge (Backchannel)'
$ws.Range("I62").Value = 'b'
$ws.Range("J62").Value = 'Acknowledge (Backchannel)'
$ws.Range("I66").Value = 'sv'
$ws.Range("J66").Value = 'Statement-opinion'
$ws.Range("I79").Value = 'sd'
$ws.Range("J79").Value = 'Statement-non-opinion'
$ws.Range("I82").Value = 'sd'
$ws.Range("J82").Value = 'Statement-non-opinion'
$ws.Range("I85").Value = 'b'
$ws.Range("J85").Value = 'Acknowledge (Backchannel)'
$ws.Range("I116").Value = 'b'
$ws.Range("J116").Value = 'Acknowledge (Backchannel)'
$ws.Range("I118").Value = 'b'
$ws.Range("J118").Value = 'Acknowledge (Backchannel)'
$ws.Range("I119").Value = 'b'
$ws.Range("J119").Value = 'Acknowledge (Backchannel)'
$ws.Range("I128").Value = 'b'
$ws.Range("J128").Value = 'Acknowledge (Backchannel)'
$ws.Range("I145").Value = '%'
$ws.Range("J145").Value = 'Uninterpretable'
$ws.Range("I157").Value = 'b'
$ws.Range("J157").Value = 'Acknowledge (Backchannel)'
$ws.Range("I159").Value = 'sd'
$ws.Range("J159").Value = 'Statement-non-opinion'
$ws.Range("I168").Value = 'b'
$ws.Range("J168").Value = 'Acknowledge (Backchannel)'
$ws.Range("I176").Value = 'b'
$ws.Range("J176").Value = 'Acknowledge (Backchannel)'
$ws.Range("I191").Value = 'sd'
$ws.Range("J191").Value = 'Statement-non-opinion'
$ws.Range("I193").Value = 'sv'
$ws.Range("J193").Value = 'Statement-opinion'
$ws.Range("I197").Value = 'b'
$ws.Range("J197").Value = 'Acknowledge (Backchannel)'
$ws.Range("I205").Value = 'b'
$ws.Range("J205").Value = 'Acknowledge (Backchannel)'
$ws.Range("I216").Value = 'b'
$ws.Range("J216").Value = 'Acknowledge (Backchannel)'
$ws.Range("I218").Value = 'b'
$ws.Range("J218").Value = 'Acknowledge (Backchannel)'
$ws.Range("I249").Value = 'sv'
$ws.Range("J249").Value = 'Statement-opinion'
$ws.Range("I267").Value = 'sv'
$ws.Range("J267").Value = 'Statement-opinion'
$ws.Range("I274").Value = 'aa'
$ws.Range("J274").Value = 'Agree/Accept'
$ws.Range("I279").Value = 'b'
$ws.Range("J279").Value = 'Acknowledge (Backchannel)'
$ws.Range("I285").Value = 'aa'
$ws.Range("J285").Value = 'Agree/Accept'
$ws.Range("I290").Value = 'sd'
$ws.Range("J290").Value = 'Statement-non-opinion'
$ws.Range("I298").Value = 'b'
$ws.Range("J298").Value = 'Acknowledge (Backchannel)'
$ws.Range("I308").Value = 'sd'
$ws.Range("J308").Value = 'Statement-non-opinion'
$ws.Range("I309").Value = 'sd'
$ws.Range("J309").Value = 'Statement-non-opinion'
$ws.Range("I314").Value = 'b'
$ws.Range("J314").Value = 'Acknowledge (Backchannel)'
$ws.Range("I320").Value = 'aa'
$ws.Range("J320").Value = 'Agree/Accept'
$ws.Range("I326").Value = 'sv'
$ws.Range("J326").Value = 'Statement-opinion'
$ws.Range("I336").Value = 'b'
$ws.Range("J336").Value = 'Acknowledge (Backchannel)'
$ws.Range("I341").Value = 'sd'
$ws.Range("J341").Value = 'Statement-non-opinion'
$ws.Range("I346").Value = 'sd'
$ws.Range("J346").Value = 'Statement-non-opinion'
$ws.Range("I351").Value = 'b'
$ws.Range("J351").Value = 'Acknowledge (Backchannel)'
$ws.Range("I353").Value = 'b'
$ws.Range("J353").Value = 'Acknowledge (Backchannel)'
